$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.884.59"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "3.902.62"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Cells.Item(5, 4).Value = "'466.80"
$ws.Range("E5").Value = "  +9.49%  "
$ws.Cells.Item(6, 4).Value = "'145.29"
$ws.Range("E6").Value = "  +6.05%  "
$ws.Cells.Item(7, 4).Value = "'0.629"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +3.77%  "
$ws.Cells.Item(11, 4).Value = "'0.0000340"
$ws.Range("E11").Value = "  +3.91%  "
$ws.Cells.Item(12, 4).Value = "'43.28"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Cells.Item(13, 4).Value = "'10.47"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "4.527.33"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Cells.Item(15, 4).Value = "'15.07"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "3.900.92"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Cells.Item(18, 4).Value = "'20.03"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").Value = "67.131.60"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Cells.Item(21, 4).Value = "'432.62"
$ws.Range("E21").Value = "  +5.18%  "
$ws.Cells.Item(22, 4).Value = "'14.72"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Cells.Item(23, 4).Value = "'3.36"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Cells.Item(24, 4).Value = "'88.59"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Cells.Item(25, 4).Value = "'38.63"
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("E26").Value = "  +6.78%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(27, 4).Value = "'10.10"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(28, 4).Value = "'5.67"
$ws.Range("E28").Value = "  +4.88%  "
$ws.Cells.Item(29, 4).Value = "'9.70"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Cells.Item(30, 4).Value = "'738.75"
$ws.Range("E30").Value = "  +4.28%  "
$ws.Cells.Item(31, 4).Value = "'13.67"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Cells.Item(34, 4).Value = "'43.79"
$ws.Range("E34").Value = "  +10.20%  "
$ws.Range("E35").Value = "  +4.23%  "
$ws.Cells.Item(36, 4).Value = "'58.00"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Cells.Item(38, 4).Value = "'3.27"
$ws.Range("E38").Value = "  +13.55%  "
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  +11.25%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).Value = "'0.0477"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("E43").Value = "  +4.79%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +6.24%  "
$ws.Range("E47").Value = "  +2.04%  "
$ws.Cells.Item(48, 4).Value = "'2.48"
$ws.Range("E48").Value = "  -5.93%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51, 4).Value = "'143.42"
$ws.Range("E51").Value = "  +0.21%  "
